$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows appended to the CELG random trade data set
$data = @(
    @(9990.1, 9959.23, 107.89, 108.22, $false, 0.31,  42613.766597222224, $true),
    @(9994.1, 9990.1,  107.17, 107.21, $false, 0.04,  42614.67386574074,  $true),
    @(9992.1, 9994.1,  107.04, 107.02, $false, -0.02, 42615.752928240741, $false)
)

$row = 6
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $row++
}
